$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 84-92 currently hold the old "CBAMT1-1".."CBAMT1-5" truckload type
# values in column D. The new type order replaces all of these with a
# single consolidated type "CBATT1" (was "CBATW1"), per the commit message.
for ($r = 84; $r -le 92; $r++) {
    $ws.Cells.Item($r, 4).Value = "CBATT1"
}

# Reset the view: scroll back to the top and select A1:D30 instead of the
# previous scrolled-down single-cell selection.
$ws.Range("A1:D30").Select()

$wb.Save()
